$p = $ppt.ActivePresentation

# --- Slide 10: "Use DP to calculate ..." textbox -----------------------
# Change "possible combination" -> "random combination" (single run,
# rewritten wholesale so the curly right-single-quote in "it's" survives
# the COM text round trip unharmed).
$rsquo = [char]0x2019
$s10 = $p.Slides.Item(10)
$shape10 = $s10.Shapes.Item("TextBox 5")
$run10 = $shape10.TextFrame.TextRange.Runs(1)
$run10.Text = "Use DP to calculate if it" + $rsquo + "s possible to sum up to cumulative value, then append random combination to matching list"

# --- Slide 11: "Issues: ..." textbox ------------------------------------
# Grow/reposition the textbox and append " and randomly matches orders".
$s11 = $p.Slides.Item(11)
$shape11 = $s11.Shapes.Item("TextBox 7")
$shape11.Top = 235.22818897637796
$shape11.Height = 167.2172440944882
$run11 = $shape11.TextFrame.TextRange.Runs(1)
$run11.Text = "Issues: mismatching price, mismatching quantity, runtime of round 5 is a bit long  and randomly matches orders"
